$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "Deadly storm slams into California's coast"
$ws.Range("B32").Value = "The storm is expected to bring deadly flooding, landsides, and power outages affecting millions."
$ws.Range("C32").Value = "https://bbc.co.uk/news/world-us-canada-64169954"
$ws.Range("D32").Value = 2
